# Update "F" (想去人数 / wanted-to-go count) figures across sheets, plus a
# couple of status/content changes, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F4").Value = 3088
$ws1.Range("F9").Value = 7265
$ws1.Range("F10").Value = 61
$ws1.Range("F13").Value = 244
$ws1.Range("F14").Value = 12
$ws1.Range("F18").Value = 1956
$ws1.Range("F19").Value = 1773
$ws1.Range("F22").Value = 1064
$ws1.Range("F24").Value = 1752
$ws1.Range("F25").Value = 1346
$ws1.Range("F28").Value = 42
$ws1.Range("F31").Value = 510
$ws1.Range("F33").Value = 2651
$ws1.Range("F34").Value = 2965
$ws1.Range("F35").Value = 2148
$ws1.Range("F36").Value = 116
$ws1.Range("F43").Value = 366
$ws1.Range("F45").Value = 232
$ws1.Range("F47").Value = 677
$ws1.Range("F49").Value = 56

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F14").Value = 90
$ws2.Range("F23").Value = 66

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F6").Value = 1816
$ws3.Range("F10").Value = 1071
# Row 13 sold out -> lowest price switches from a number to "已售罄"
$ws3.Range("F13").Value = 1797
$ws3.Range("G13").Value = "已售罄"
$ws3.Range("F14").Value = 8047

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# Row 10 becomes a brand-new event (MADEBYBILIBILI signing), and the event
# that used to sit in row 10 ("剑网3 x HAPPY ZOO") moves down into row 11,
# replacing the event that used to be there ("东方明珠" pop-up, which is
# dropped from this aggregate sheet).
$ws4.Range("C10").Value = "上海·MADEBYBILIBILI高能中心徐汇万科中心站·高能国漫签售会"
$ws4.Range("D10").Value = "沪闵路9191号 徐汇万科广场"
$ws4.Range("E10").Value = "2024.08.16 10:30-08.25 12:00"
$ws4.Range("F10").Value = 142
$ws4.Range("G10").Value = 1
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=90589"
$ws4.Range("I10").Value = "//i0.hdslb.com/bfs/openplatform/202408/qYopb72A1723451211534.png"

# Force this as literal text (it already looks like a date, and the
# worksheet stores it as text) so Excel doesn't auto-convert it into a
# date serial number.
$ws4.Range("B11").NumberFormat = "@"
$ws4.Range("B11").Value = "2024-08-16"
$ws4.Range("C11").Value = "上海·剑网3×HAPPY ZOO 剑网3十五周年主题咖啡厅"
$ws4.Range("D11").Value = "南京东路340号百联zx创趣场四楼05号 HAPPY ZOO"
$ws4.Range("E11").Value = "2024.08.16 00:00-10.07 23:59"
$ws4.Range("F11").Value = 398
$ws4.Range("G11").Value = 10
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=90305"
$ws4.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202408/QzCwiYge1722838646403.png"

$ws4.Range("F3").Value = 3088
$ws4.Range("F4").Value = 1816
$ws4.Range("F7").Value = 7265
$ws4.Range("F9").Value = 61
$ws4.Range("F13").Value = 12
$ws4.Range("F18").Value = 1956
$ws4.Range("F20").Value = 1064
$ws4.Range("F22").Value = 1752
$ws4.Range("F23").Value = 1346
$ws4.Range("F27").Value = 42
$ws4.Range("F29").Value = 90
$ws4.Range("F31").Value = 510
$ws4.Range("F34").Value = 2651
$ws4.Range("F35").Value = 2965
$ws4.Range("F36").Value = 2148
$ws4.Range("F37").Value = 116
$ws4.Range("F43").Value = 366
$ws4.Range("F45").Value = 66
$ws4.Range("F46").Value = 232
